$d = $word.ActiveDocument
$s = $d.Styles.Item("Kop6")
Write-Output "ListTemplate: $($s.ListTemplate)"
Write-Output "ListLevelNumber: $($s.ListLevelNumber)"
try {
  $s.ListLevelNumber = 0
  Write-Output "set ListLevelNumber=0 ok"
} catch {
  Write-Output "ERR: $_"
}
